# Update column G ("K" - strikeouts) values for rows 2-24 on the active sheet.
# This mirrors a regen of save_data where Strike# was replaced by K and the
# per-game strikeout totals (K) were recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 4
    6  = 2
    7  = 1
    8  = 2
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 3
    14 = 2
    15 = 1
    16 = 0
    17 = 2
    18 = 2
    19 = 0
    20 = 0
    21 = 0
    22 = 2
    23 = 1
    24 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
